$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 329747
$ws.Range("D2").Value = 419937615
$ws.Range("C4").Value = 332
$ws.Range("D4").Value = 474987
$ws.Range("C8").Value = 877
$ws.Range("D8").Value = 1289899
$ws.Range("C10").Value = 119395
$ws.Range("D10").Value = 174938243
$ws.Range("C12").Value = 61237
$ws.Range("D12").Value = 88372499
$ws.Range("C16").Value = 4051
$ws.Range("D16").Value = 5751137
$ws.Range("C20").Value = 7153
$ws.Range("D20").Value = 9988435
$ws.Range("C22").Value = 79301
$ws.Range("D22").Value = 98769465
$ws.Range("C28").Value = 33001
$ws.Range("D28").Value = 48301181
$ws.Range("C30").Value = 11761
$ws.Range("D30").Value = 16916106
$ws.Range("C35").Value = 1966
$ws.Range("D35").Value = 2775092
$ws.Range("C36").Value = 99303
$ws.Range("D36").Value = 124833671
$ws.Range("C44").Value = 45092
$ws.Range("D44").Value = 66075533
$ws.Range("C46").Value = 9388
$ws.Range("D46").Value = 13463729
$ws.Range("C51").Value = 2502
$ws.Range("D51").Value = 3497418
$ws.Range("C52").Value = 70656
$ws.Range("D52").Value = 88586409
$ws.Range("C59").Value = 28707
$ws.Range("D59").Value = 42102120
$ws.Range("C62").Value = 11429
$ws.Range("D62").Value = 16523918
$ws.Range("C64").Value = 1374
$ws.Range("D64").Value = 1920997
$ws.Range("C68").Value = 1588
$ws.Range("D68").Value = 2227351
$ws.Range("C70").Value = 20946
$ws.Range("D70").Value = 27428748
$ws.Range("C74").Value = 7741
$ws.Range("D74").Value = 11337300
$ws.Range("C76").Value = 5236
$ws.Range("D76").Value = 7604494
$ws.Range("C79").Value = 144262
$ws.Range("D79").Value = 179743613
$ws.Range("C83").Value = 442
$ws.Range("D83").Value = 645824
$ws.Range("C85").Value = 64752
$ws.Range("D85").Value = 94894265
$ws.Range("C88").Value = 30454
$ws.Range("D88").Value = 44054104
$ws.Range("C90").Value = 2772
$ws.Range("D90").Value = 3990652
$ws.Range("C91").Value = 3011
$ws.Range("D91").Value = 4255668
$ws.Range("C92").Value = 34624
$ws.Range("D92").Value = 46951793
$ws.Range("C96").Value = 8394
$ws.Range("D96").Value = 12339821
$ws.Range("C98").Value = 7779
$ws.Range("D98").Value = 11290347
$ws.Range("C100").Value = 554
$ws.Range("D100").Value = 786156
$ws.Range("C101").Value = 524
$ws.Range("D101").Value = 756550
$ws.Range("C102").Value = 11642
$ws.Range("D102").Value = 18764262
$ws.Range("C104").Value = 2827
$ws.Range("D104").Value = 4893970
$ws.Range("C106").Value = 3891
$ws.Range("D106").Value = 6764869
$ws.Range("C108").Value = 171
$ws.Range("D108").Value = 294045
$ws.Range("C109").Value = 225
$ws.Range("D109").Value = 364530
$ws.Range("C110").Value = 145328
$ws.Range("D110").Value = 179733269
$ws.Range("C114").Value = 962
$ws.Range("D114").Value = 1410815
$ws.Range("C116").Value = 53777
$ws.Range("D116").Value = 78809549
$ws.Range("C117").Value = 89
$ws.Range("D117").Value = 131959
$ws.Range("C118").Value = 28023
$ws.Range("D118").Value = 40600436
$ws.Range("C122").Value = 2412
$ws.Range("D122").Value = 3392994
$ws.Range("C124").Value = 540898
$ws.Range("D124").Value = 714802948
$ws.Range("C125").Value = 94
$ws.Range("D125").Value = 125071
$ws.Range("C126").Value = 223
$ws.Range("D126").Value = 328509
$ws.Range("C129").Value = 1412
$ws.Range("D129").Value = 2092714
$ws.Range("C131").Value = 215029
$ws.Range("D131").Value = 316067834
$ws.Range("C132").Value = 425
$ws.Range("D132").Value = 634210
$ws.Range("C134").Value = 193252
$ws.Range("D134").Value = 281033776
$ws.Range("C136").Value = 35
$ws.Range("D136").Value = 51332
$ws.Range("C137").Value = 2890
$ws.Range("D137").Value = 4057572
$ws.Range("C140").Value = 6842
$ws.Range("D140").Value = 9654555
$ws.Range("C143").Value = 46189
$ws.Range("D143").Value = 61642790
$ws.Range("C149").Value = 14434
$ws.Range("D149").Value = 21156239
$ws.Range("C150").Value = 3873
$ws.Range("D150").Value = 5585482
$ws.Range("C155").Value = 416
$ws.Range("D155").Value = 586813
$ws.Range("C156").Value = 18205
$ws.Range("D156").Value = 24065567
$ws.Range("C160").Value = 7427
$ws.Range("D160").Value = 10808593
$ws.Range("C162").Value = 5196
$ws.Range("D162").Value = 7479538
$ws.Range("C167").Value = 21301
$ws.Range("D167").Value = 37579924
$ws.Range("C168").Value = 2271
$ws.Range("D168").Value = 3998571
$ws.Range("C169").Value = 297
$ws.Range("D169").Value = 510089
$ws.Range("C172").Value = 120
$ws.Range("D172").Value = 218949
$ws.Range("C173").Value = 89871
$ws.Range("D173").Value = 112236808
$ws.Range("C180").Value = 34508
$ws.Range("D180").Value = 50597758
$ws.Range("C182").Value = 13384
$ws.Range("D182").Value = 19337488
$ws.Range("C184").Value = 1268
$ws.Range("D184").Value = 1774527
$ws.Range("C186").Value = 1753
$ws.Range("D186").Value = 2460929
$ws.Range("C188").Value = 244090
$ws.Range("D188").Value = 303223086
$ws.Range("C196").Value = 88137
$ws.Range("D196").Value = 129178702
$ws.Range("C199").Value = 33865
$ws.Range("D199").Value = 48754132
$ws.Range("C202").Value = 5193
$ws.Range("D202").Value = 7393505
$ws.Range("C205").Value = 5202
$ws.Range("D205").Value = 7205801
$ws.Range("C208").Value = 270566
$ws.Range("D208").Value = 334738435
$ws.Range("C215").Value = 626
$ws.Range("D215").Value = 911878
$ws.Range("C217").Value = 97034
$ws.Range("D217").Value = 141949861
$ws.Range("C220").Value = 52979
$ws.Range("D220").Value = 76565846
$ws.Range("C223").Value = 4739
$ws.Range("D223").Value = 6650835
$ws.Range("C226").Value = 6193
$ws.Range("D226").Value = 8582999
$ws.Range("C229").Value = 109106
$ws.Range("D229").Value = 136368209
$ws.Range("C231").Value = 77
$ws.Range("D231").Value = 110513
$ws.Range("C236").Value = 50420
$ws.Range("D236").Value = 73858454
$ws.Range("C237").Value = 39
$ws.Range("D237").Value = 56211
$ws.Range("C238").Value = 12851
$ws.Range("D238").Value = 18485535
$ws.Range("C240").Value = 1911
$ws.Range("D240").Value = 2740382
$ws.Range("C242").Value = 2671
$ws.Range("D242").Value = 3741506
$ws.Range("C243").Value = 265132
$ws.Range("D243").Value = 334739703
$ws.Range("C249").Value = 846
$ws.Range("D249").Value = 1242404
$ws.Range("C250").Value = 11
$ws.Range("D250").Value = 16500
$ws.Range("C251").Value = 97909
$ws.Range("D251").Value = 143454141
$ws.Range("C252").Value = 78
$ws.Range("D252").Value = 114161
$ws.Range("C254").Value = 67267
$ws.Range("D254").Value = 97506155
$ws.Range("C256").Value = 2452
$ws.Range("D256").Value = 3458724
$ws.Range("C259").Value = 4944
$ws.Range("D259").Value = 6943190
